$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.572.75"
$ws.Range("E2").Value = "  -5.55%  "
$ws.Range("D3").Value = "3.056.88"
$ws.Range("E3").Value = "  -6.00%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.66%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.051.82"
$ws.Range("E8").Value = "  -5.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.39%  "
$ws.Range("E10").Value = "  -6.17%  "
$ws.Range("E11").Value = "  -13.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000220"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.41%  "
$ws.Range("D15").Value = "3.505.87"
$ws.Range("E15").Value = "  -7.13%  "
$ws.Range("D16").Value = "62.597.16"
$ws.Range("E16").Value = "  -5.61%  "
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "3.060.43"
$ws.Range("E18").Value = "  -5.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -13.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.58%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -14.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.32%  "
$ws.Range("E32").Value = "  -7.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -12.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "489.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -13.19%  "
$ws.Range("E36").Value = "  -6.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.31%  "
$ws.Range("D38").Value = "3.135.83"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.58%  "
$ws.Range("E40").Value = "  -7.89%  "
$ws.Range("E41").Value = "  -10.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -16.10%  "
$ws.Range("E44").Value = "  -10.57%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -12.09%  "
$ws.Range("E47").Value = "  -7.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.107"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.10%  "
$ws.Range("D50").Value = "0.0₃0505"
$ws.Range("E50").Value = "  -9.80%  "
$ws.Range("E51").Value = "  -9.40%  "
